# Informe-03-030026-A-TC-TP.xlsx — metadata sheet rework
#
# #8  Mejorar la generacion de SKOS Concept Schemes
# #16 Incluir descripciones para algunas medidas en los DSDs
# #17 Referenciada codelist que luego no tiene valores
# #19 Anadir propiedad en el DSD que identifique el ambito territorial aplicable
# #20 Generacion erronea de medidas en 01-080101-010105TC
#
# Row 1 = human-readable column labels (now capitalised / re-worded)
# Row 2 = measure/dimension identifier backing each column ("null" if n/a)
# Row 3 = "medida" (measure) vs "dim" (dimension) classifier
# Row 4 = datatype (xsd:string / xsd:int) for measures, URI-* codelist for dims

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: column labels -------------------------------------------------
$ws.Range("A1").Value = "Continente"
$ws.Range("B1").Value = "Nacionalidad, nombre"
$ws.Range("C1").Value = "Personas"
$ws.Range("D1").Value = "Área"
$ws.Range("E1").Value = "Nacionalidad, código"
$ws.Range("F1").Value = "Area nacionalidad, código"
$ws.Range("G1").Value = "Comarca nombre"
$ws.Range("H1").Value = "Comarca código"
$ws.Range("I1").Value = "Provincia código"
$ws.Range("J1").Value = "Aragón"
$ws.Range("K1").Value = "Provincia nombre"

# --- Row 2: measure / dimension identifier --------------------------------
$ws.Range("B2").Value = "iaest-measure:nacionalidad-nombre"
$ws.Range("C2").Value = "iaest-measure:personas"
$ws.Range("D2").Value = "iaest-measure:area"
$ws.Range("F2").Value = "null"
$ws.Range("G2").Value = "sdmx-dimension:refArea"
$ws.Range("H2").Value = "null"
$ws.Range("J2").Value = "sdmx-dimension:refArea"

# --- Row 3: medida / dim classifier ---------------------------------------
$ws.Range("C3").Value = "medida"
$ws.Range("F3").Value = "null"
$ws.Range("G3").Value = "dim"
$ws.Range("H3").Value = "null"
$ws.Range("J3").Value = "dim"

# --- Row 4: datatype / codelist URI ---------------------------------------
$ws.Range("A4").Value = "xsd:string"
$ws.Range("B4").Value = "xsd:string"
$ws.Range("C4").Value = "xsd:int"
$ws.Range("D4").Value = "xsd:string"
$ws.Range("F4").Value = "null"
$ws.Range("G4").Value = "URI-comarca"
$ws.Range("H4").Value = "null"
$ws.Range("J4").Value = "URI-Comunidad"
$ws.Range("K4").Value = "URI-Provincia"
